$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 71430000
$ws.Range("J103").Value = 83334830
$ws.Range("L103").Value = 250004490
$ws.Range("N103").Value = -250005662
$ws.Range("H107").Value = 2241
$ws.Range("J107").Value = 2982.3333
$ws.Range("L107").Value = 2982.3333
$ws.Range("N107").Value = -6822.3333
$ws.Range("H123").Value = 99999
$ws.Range("J123").Value = 99999
$ws.Range("L123").Value = 99999
$ws.Range("N123").Value = -109799
$ws.Range("H129").Value = 2613
$ws.Range("I129").Value = 967.875
$ws.Range("J129").Value = 9193.5
$ws.Range("K129").Value = 2903.625
$ws.Range("L129").Value = 27580.5
$ws.Range("M129").Value = 2096.375
$ws.Range("N129").Value = -37580.5
$ws.Range("H132").Value = 4841.7095
$ws.Range("I132").Value = 2190.762
$ws.Range("J132").Value = 10408.7
$ws.Range("K132").Value = 6572.286
$ws.Range("L132").Value = 31226.1
$ws.Range("M132").Value = -4042.286
$ws.Range("N132").Value = -36286.10000000001
$ws.Range("H138").Value = 3792.5356
$ws.Range("I138").Value = 1520.1923
$ws.Range("J138").Value = 33333
$ws.Range("K138").Value = 4560.5769
$ws.Range("L138").Value = 99999
$ws.Range("M138").Value = 579.4231
$ws.Range("N138").Value = -110279

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4165.8237
$ws.Range("I32").Value = 4243.5557
$ws.Range("K32").Value = 4243.5557
$ws.Range("M32").Value = -3956.5557
$ws.Range("H45").Value = 1498.375
$ws.Range("I45").Value = 1426.8572
$ws.Range("J45").Value = 1999
$ws.Range("K45").Value = 1426.8572
$ws.Range("L45").Value = 1999
$ws.Range("M45").Value = -1049.8572
$ws.Range("N45").Value = -2753
$ws.Range("H61").Value = 26253674
$ws.Range("I61").Value = 40004800
$ws.Range("J61").Value = 3335132.8
$ws.Range("K61").Value = 40004800
$ws.Range("L61").Value = 3335132.8
$ws.Range("M61").Value = -40004588
$ws.Range("N61").Value = -3335556.8
$ws.Range("H74").Value = 782926.9
$ws.Range("I74").Value = 893844.5600000001
$ws.Range("J74").Value = 6503
$ws.Range("K74").Value = 893844.5600000001
$ws.Range("L74").Value = 6503
$ws.Range("M74").Value = -892970.5600000001
$ws.Range("N74").Value = -8251
$ws.Range("H77").Value = 782926.9
$ws.Range("I77").Value = 893844.5600000001
$ws.Range("J77").Value = 6503
$ws.Range("K77").Value = 4469222.800000001
$ws.Range("L77").Value = 32515
$ws.Range("M77").Value = -4464854.800000001
$ws.Range("N77").Value = -41251
$ws.Range("H110").Value = 1783
$ws.Range("I110").Value = 764.1429000000001
$ws.Range("J110").Value = 5349
$ws.Range("K110").Value = 764.1429000000001
$ws.Range("L110").Value = 5349
$ws.Range("M110").Value = 1280.8571
$ws.Range("N110").Value = -9439
$ws.Range("H132").Value = 1757510.9
$ws.Range("I132").Value = 3063.6086
$ws.Range("K132").Value = 9190.825800000001
$ws.Range("M132").Value = -6660.825800000001
$ws.Range("H136").Value = 26253674
$ws.Range("I136").Value = 40004800
$ws.Range("J136").Value = 3335132.8
$ws.Range("K136").Value = 120014400
$ws.Range("L136").Value = 10005398.4
$ws.Range("M136").Value = -120011850
$ws.Range("N136").Value = -10010498.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2151.2
$ws.Range("I94").Value = 1671
$ws.Range("K94").Value = 1671
$ws.Range("M94").Value = -1220
$ws.Range("H134").Value = 3450768.2
$ws.Range("I134").Value = 2416.1667
$ws.Range("K134").Value = 7248.500100000001
$ws.Range("M134").Value = -4713.500100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 77.666664
$ws.Range("I7").Value = 77.666664
$ws.Range("K7").Value = 77.666664
$ws.Range("M7").Value = 35.333336
$ws.Range("H132").Value = 3374.5
$ws.Range("I132").Value = 2999.8333
$ws.Range("J132").Value = 4498.5
$ws.Range("K132").Value = 8999.499899999999
$ws.Range("L132").Value = 13495.5
$ws.Range("M132").Value = -6469.499899999999
$ws.Range("N132").Value = -18555.5
$ws.Range("H134").Value = 2025.3077
$ws.Range("I134").Value = 2054.6365
$ws.Range("J134").Value = 1864
$ws.Range("K134").Value = 6163.9095
$ws.Range("L134").Value = 5592
$ws.Range("M134").Value = -3628.9095
$ws.Range("N134").Value = -10662

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2616.9546
$ws.Range("I80").Value = 1705.5
$ws.Range("K80").Value = 1705.5
$ws.Range("M80").Value = -707.5
$ws.Range("H83").Value = 2616.9546
$ws.Range("I83").Value = 1705.5
$ws.Range("K83").Value = 8527.5
$ws.Range("M83").Value = -3535.5
$ws.Range("H102").Value = 2625.8572
$ws.Range("I102").Value = 2651
$ws.Range("K102").Value = 2651
$ws.Range("M102").Value = -1029
$ws.Range("H113").Value = 928469.25
$ws.Range("I113").Value = 2837.2307
$ws.Range("J113").Value = 2647500.2
$ws.Range("K113").Value = 2837.2307
$ws.Range("L113").Value = 2647500.2
$ws.Range("M113").Value = -667.2307000000001
$ws.Range("N113").Value = -2651840.2
$ws.Range("H122").Value = 4852.3
$ws.Range("I122").Value = 5579.5
$ws.Range("K122").Value = 16738.5
$ws.Range("M122").Value = -14288.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5611.2
$ws.Range("I7").Value = 5611.2
$ws.Range("K7").Value = 5611.2
$ws.Range("M7").Value = -5499.2
$ws.Range("H16").Value = 3636.318
$ws.Range("J16").Value = 7907.6665
$ws.Range("L16").Value = 7907.6665
$ws.Range("N16").Value = -8247.666499999999
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H46").Value = 6299.6
$ws.Range("J46").Value = 7500
$ws.Range("L46").Value = 7500
$ws.Range("N46").Value = -7876
$ws.Range("H100").Value = 19253036
$ws.Range("I100").Value = 2441.375
$ws.Range("J100").Value = 50053988
$ws.Range("K100").Value = 2441.375
$ws.Range("L100").Value = 50053988
$ws.Range("M100").Value = -1900.375
$ws.Range("N100").Value = -50055070
$ws.Range("H126").Value = 5611.2
$ws.Range("I126").Value = 5611.2
$ws.Range("K126").Value = 16833.6
$ws.Range("M126").Value = -14363.6
$ws.Range("H132").Value = 4899.5386
$ws.Range("I132").Value = 3144.1428
$ws.Range("K132").Value = 9432.428400000001
$ws.Range("M132").Value = -6902.428400000001
$ws.Range("H136").Value = 3384.7273
$ws.Range("I136").Value = 2882.2307
$ws.Range("J136").Value = 4110.5557
$ws.Range("K136").Value = 8646.6921
$ws.Range("L136").Value = 12331.6671
$ws.Range("M136").Value = -6096.6921
$ws.Range("N136").Value = -17431.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2471.5789
$ws.Range("I107").Value = 1572.1765
$ws.Range("J107").Value = 3199.6667
$ws.Range("K107").Value = 4716.529500000001
$ws.Range("L107").Value = 9599.000100000001
$ws.Range("M107").Value = -2796.529500000001
$ws.Range("N107").Value = -13439.0001
$ws.Range("H113").Value = 997.5
$ws.Range("I113").Value = 996.6667
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 2990.0001
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -820.0001000000002
$ws.Range("N113").Value = -7340
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H126").Value = 12363.4
$ws.Range("I126").Value = 12999
$ws.Range("J126").Value = 11410
$ws.Range("K126").Value = 38997
$ws.Range("L126").Value = 34230
$ws.Range("M126").Value = -36527
$ws.Range("N126").Value = -39170
$ws.Range("H132").Value = 501707.16
$ws.Range("I132").Value = 1714.3529
$ws.Range("J132").Value = 3334999.8
$ws.Range("K132").Value = 5143.0587
$ws.Range("L132").Value = 10004999.4
$ws.Range("M132").Value = -2613.0587
$ws.Range("N132").Value = -10010059.4
$ws.Range("H140").Value = 96997
$ws.Range("J140").Value = 96997
$ws.Range("L140").Value = 96997
$ws.Range("N140").Value = -107357
